# KWD-loginGmail.xlsx update: "updated to latest olo"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Workbook window was nudged/resized slightly by the author's Excel
#     session before saving. ---
$win = $wb.Windows.Item(1)
$win.Left = 240
$win.Top = 165
$win.Width = 14805
$win.Height = 7950

# --- Remove the old trailing rows (12-23); the refreshed keyword sheet only
#     needs 10 data rows (+ header) now. Deleting the whole rows (not just
#     clearing) shrinks the used range down to A1:D11 like the new sheet. ---
$ws.Range("A12:D23").EntireRow.Delete() | Out-Null

# --- Clear the stale date-format on the old C7 cell so it falls back to the
#     plain text style used by the rest of the data column. ---
$ws.Cells.Item(7, 3).NumberFormat = "@"

# --- Seed the brand-new keyword/value strings in the order the refreshed
#     test case introduces them, then fill in the rest of the grid below. ---
$ws.Cells.Item(3, 3).Value = "pavan.mailme"
$ws.Cells.Item(1, 1).Value = "Command"
$ws.Cells.Item(1, 2).Value = "Target"
$ws.Cells.Item(4, 3).Value = "Bing"
$ws.Cells.Item(10, 1).Value = "AssertTitle"
$ws.Cells.Item(4, 1).Value = "VerifyTitle"

# --- Rewrite the 11 remaining rows with the refreshed keyword data. ---

# Header row
$ws.Cells.Item(1, 1).Value = "Command"
$ws.Cells.Item(1, 2).Value = "Target"
$ws.Cells.Item(1, 3).Value = "Value"
$ws.Cells.Item(1, 4).Value = "Options"

# Row 2
$ws.Cells.Item(2, 1).Value = "IfElementPresent"
$ws.Cells.Item(2, 2).Value = "login.username"
$ws.Cells.Item(2, 3).Value = ""
$ws.Cells.Item(2, 4).Value = ""

# Row 3
$ws.Cells.Item(3, 1).Value = "Type"
$ws.Cells.Item(3, 2).Value = "login.username"
$ws.Cells.Item(3, 3).Value = "pavan.mailme"
$ws.Cells.Item(3, 4).Value = ""

# Row 4
$ws.Cells.Item(4, 1).Value = "VerifyTitle"
$ws.Cells.Item(4, 2).Value = ""
$ws.Cells.Item(4, 3).Value = "Bing"
$ws.Cells.Item(4, 4).Value = ""

# Row 5
$ws.Cells.Item(5, 1).Value = "Else"
$ws.Cells.Item(5, 2).Value = ""
$ws.Cells.Item(5, 3).Value = ""
$ws.Cells.Item(5, 4).Value = ""

# Row 6
$ws.Cells.Item(6, 1).Value = "VerifyElementNotPresent"
$ws.Cells.Item(6, 2).Value = "login.password"
$ws.Cells.Item(6, 3).Value = ""
$ws.Cells.Item(6, 4).Value = ""

# Row 7
$ws.Cells.Item(7, 1).Value = "VerifyElementNotPresent"
$ws.Cells.Item(7, 2).Value = "login.username"
$ws.Cells.Item(7, 3).Value = ""
$ws.Cells.Item(7, 4).Value = ""

# Row 8
$ws.Cells.Item(8, 1).Value = "EndIf"
$ws.Cells.Item(8, 2).Value = ""
$ws.Cells.Item(8, 3).Value = ""
$ws.Cells.Item(8, 4).Value = ""

# Row 9
$ws.Cells.Item(9, 1).Value = "Type"
$ws.Cells.Item(9, 2).Value = "login.password"
$ws.Cells.Item(9, 3).Value = "pavan.mailme"
$ws.Cells.Item(9, 4).Value = ""

# Row 10
$ws.Cells.Item(10, 1).Value = "AssertTitle"
$ws.Cells.Item(10, 2).Value = ""
$ws.Cells.Item(10, 3).Value = "Bing"
$ws.Cells.Item(10, 4).Value = ""

# Row 11
$ws.Cells.Item(11, 1).Value = "VerifyElementPresent"
$ws.Cells.Item(11, 2).Value = "login.password"
$ws.Cells.Item(11, 3).Value = ""
$ws.Cells.Item(11, 4).Value = ""

# --- Move the active selection to C10, matching the refreshed view state. ---
$ws.Range("C10").Select() | Out-Null
